$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume data updated by the scheduled scraper run.
# D/E columns store text (not numbers) - some D values look numeric (e.g. "0.9980")
# so we force text format first to avoid Excel auto-converting / trimming them.

$ws.Range('D2').Value = '29.044.17'
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').Value = '1.831.94'
$ws.Range('E3').Value = '  +0.47%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9980'
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.56'
$ws.Range('E5').Value = '  -0.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6181'
$ws.Range('E6').Value = '  -1.63%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9996'
$ws.Range('E7').Value = '  +0.34%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07443'
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2933'
$ws.Range('E9').Value = '  +0.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.04'
$ws.Range('E10').Value = '  +0.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07664'
$ws.Range('E11').Value = '  -0.06%  '
$ws.Range('D12').Value = '1.850.85'
$ws.Range('E12').Value = '  +1.48%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.001'
$ws.Range('E13').Value = '  +0.52%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6727'
$ws.Range('E14').Value = '  +1.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '82.82'
$ws.Range('E15').Value = '  -0.01%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000009149'
$ws.Range('E16').Value = '  -5.35%  '
$ws.Range('E17').Value = '  -1.79%  '
$ws.Range('D18').Value = '29.068.82'
$ws.Range('E18').Value = '  +0.51%  '
$ws.Range('D19').Value = '2.087.70'
$ws.Range('E19').Value = '  +1.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '239.81'
$ws.Range('E20').Value = '  +6.73%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.69'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9997'
$ws.Range('E22').Value = '  +0.52%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.199'
$ws.Range('E23').Value = '  +1.33%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9996'
$ws.Range('E24').Value = '  +0.26%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '158.57'
$ws.Range('E25').Value = '  -0.89%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1412'
$ws.Range('E26').Value = '  +0.56%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.501'
$ws.Range('E27').Value = '  +0.32%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.86'
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.497'
$ws.Range('E29').Value = '  +0.36%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05629'
$ws.Range('E30').Value = '  +3.53%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.111'
$ws.Range('E31').Value = '  +1.60%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.127'
$ws.Range('E32').Value = '  +0.45%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.205'
$ws.Range('E33').Value = '  +0.73%  '
$ws.Range('E34').Value = '  -0.33%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7405'
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.141'
$ws.Range('E36').Value = '  +0.93%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.653'
$ws.Range('E37').Value = '  +1.81%  '
$ws.Range('E38').Value = '  +1.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01784'
$ws.Range('E39').Value = '  +0.66%  '
$ws.Range('D40').Value = '1.210.68'
$ws.Range('E40').Value = '  -2.20%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.396'
$ws.Range('E41').Value = '  -3.74%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8968'
$ws.Range('E42').Value = '  +0.17%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9981'
$ws.Range('E43').Value = '  +0.36%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.33'
$ws.Range('E44').Value = '  +0.17%  '
$ws.Range('E45').Value = '  +0.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '65.41'
$ws.Range('E46').Value = '  +0.99%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5081'
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('B48').Value = 'TheSandbox'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4058'
$ws.Range('E48').Value = '  +0.59%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.149'
$ws.Range('E49').Value = '  +2.26%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00000000117'
$ws.Range('E50').Value = '  -3.76%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05811'
$ws.Range('E51').Value = '  +0.43%  '
